# Insert a new data row just before the current row 144, shifting all
# subsequent rows (old 144-171) down by one (new 145-172), and populate
# the newly inserted row with its own values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 144 (pushes existing row 144 -> 145, etc.)
$ws.Rows.Item(144).Insert()

# Populate the newly inserted row 144 with the new record's data
$ws.Range("A144").Value = 5
$ws.Range("B144").Value = "Macroferia Regional de Talca"
$ws.Range("C144").Value = "Maule"
$ws.Range("D144").Value2 = 44946
$ws.Range("E144").Value = 7
$ws.Range("F144").Value = 100112030
$ws.Range("G144").Value = "Poroto granado"
$ws.Range("H144").Value = "Sin especificar"
$ws.Range("I144").Value = "Primera"
$ws.Range("J144").Value = 500
$ws.Range("K144").Value = 40000
$ws.Range("L144").Value = 40000
$ws.Range("M144").Value = 40000
$ws.Range("N144").Value = "$/saco 25 kilos"
$ws.Range("O144").Value = "Región del Maule"
$ws.Range("P144").Value = 1600
$ws.Range("Q144").Value = 25
$ws.Range("R144").Value = "Hortaliza"

# Match the date number formatting used by the rest of column D
$ws.Range("D144").NumberFormat = $ws.Range("D145").NumberFormat()
